# Applies the "used_for / category" row restructuring + related text edits
# to the "Examples & Info" sheet of the sequencing spreadsheet template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples & Info")

# ------------------------------------------------------------------
# 1. N2: clarify the SEX column description
# ------------------------------------------------------------------
$ws.Range("N2").Value = "Biological sex of the sampled individual."

# ------------------------------------------------------------------
# 2. Row 4 ("regime" -> "used_for"): relabel and fill in the
#    previously-blank "used_for" cells with their owning group.
# ------------------------------------------------------------------
$ws.Range("A4").Value = "used_for"

$ws.Range("R4:W4").Value = "Odomlab"
$ws.Range("X4").Value = "Odomlab, GUIDE, OTP"
$ws.Range("Y4:AK4").Value = "Odomlab"
$ws.Range("AR4").Value = "GUIDE, Odomlab"
$ws.Range("AS4:BB4").Value = "Odomlab"
$ws.Range("BD4:BJ4").Value = "Odomlab"

# ------------------------------------------------------------------
# 3. Insert a new row 6 ("category"), pushing the old "regex" row
#    (row 6) down to row 7. The new "category" row reuses the
#    sample/experiment/sequencing/other grouping that used to live
#    in row 5 ("order"), which is now cleared down to just B5.
# ------------------------------------------------------------------
$ws.Rows.Item(6).Insert()

$ws.Range("A5:BJ5").Copy()
$ws.Range("A6:BJ6").PasteSpecial(-4122)
$ws.Range("A5:BJ5").Copy()
$ws.Range("A6:BJ6").PasteSpecial(-4163)
$ws.Rows.Item(6).RowHeight = 25

$ws.Range("A6").Value = "category"
$ws.Range("B6").Value = ""

$ws.Range("A5").Value = ""
$ws.Range("C5:BJ5").ClearContents()
